# Add a run containing the text "R" (Arial, 12pt / minorBidi theme font)
# to the empty bulleted paragraph that immediately follows
#   "There are two circuit design schemes we confirmed and it is hard to
#    choose which one to use."
# in the Week 1 "Problem, issues and concerns:" list.

$d = $word.ActiveDocument

# Locate the unique anchor sentence that precedes the (currently empty)
# target paragraph.
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "There are two circuit design schemes we confirmed and it is hard to choose which one to use.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the anchor paragraph text."
}

# $anchor now spans the matched sentence; $anchor.End sits right before the
# paragraph mark that ends that paragraph, so the next paragraph (the empty
# bullet we need to fill in) starts at $anchor.End + 1.
$targetStart = $anchor.End + 1
$target = $d.Range($targetStart, $targetStart)

# Sanity-check that the target paragraph is indeed empty (just the paragraph
# mark) before we touch it. Note: a zero-length (collapsed) Range's .Text can
# report stale data in this host, so probe a 1-char range instead.
$probe = $d.Range($targetStart, $targetStart + 1)
if ($probe.Text -ne "" -and $probe.Text -ne "`r") {
    throw "Target paragraph was not empty; aborting to avoid clobbering content."
}

# Insert a fully-specified run (text + run formatting) as literal WordprocessingML
# so the run properties match exactly: Arial / minorBidi theme font, 12pt (sz=24/szCs=24).
$runXml = '<?xml version="1.0"?>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' +
            '<w:r>' +
                '<w:rPr>' +
                    '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/>' +
                    '<w:sz w:val="24"/>' +
                    '<w:szCs w:val="24"/>' +
                '</w:rPr>' +
                '<w:t>R</w:t>' +
            '</w:r>' +
        '</w:p></w:body>' +
    '</w:document>'

$target.InsertXML($runXml)
